$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 146, shifting rows 146:206 down to 147:207
$ws.Rows.Item(146).Insert()

# Populate the newly inserted row 146 with the new weekly record
$ws.Cells.Item(146, 1).Value = 3
$ws.Cells.Item(146, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(146, 3).Value = "Coquimbo"
$ws.Cells.Item(146, 4).Value = 44466
$ws.Cells.Item(146, 5).Value = 5
$ws.Cells.Item(146, 6).Value = 100112040
$ws.Cells.Item(146, 7).Value = "Cilantro"
$ws.Cells.Item(146, 8).Value = "Sin especificar"
$ws.Cells.Item(146, 9).Value = "Primera"
$ws.Cells.Item(146, 10).Value = 160
$ws.Cells.Item(146, 11).Value = 2500
$ws.Cells.Item(146, 12).Value = 2500
$ws.Cells.Item(146, 13).Value = 2500
$ws.Cells.Item(146, 14).Value = "$/docena de atados (3 kilos)"
$ws.Cells.Item(146, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(146, 16).Value = 833
$ws.Cells.Item(146, 17).Value = 3
$ws.Cells.Item(146, 18).Value = "Hortaliza"
